# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F column (time_taken) timestamps on the "data" sheet ---
$timestamps = @{
    2 = "2021-10-05 14:19:33.554152"
    3 = "2021-10-05 14:19:33.554159"
    4 = "2021-10-05 14:19:33.554162"
    5 = "2021-10-05 14:19:33.554165"
    6 = "2021-10-05 14:19:33.554168"
    7 = "2021-10-05 14:19:33.554170"
    8 = "2021-10-05 14:19:33.554173"
    9 = "2021-10-05 14:19:33.554176"
    10 = "2021-10-05 14:19:33.554179"
    11 = "2021-10-05 14:19:33.554181"
    12 = "2021-10-05 14:19:33.554184"
    13 = "2021-10-05 14:19:33.554186"
    14 = "2021-10-05 14:19:33.554189"
    15 = "2021-10-05 14:19:33.554191"
    16 = "2021-10-05 14:19:33.554194"
    17 = "2021-10-05 14:19:33.554196"
    18 = "2021-10-05 14:19:33.554200"
    19 = "2021-10-05 14:19:33.554202"
    20 = "2021-10-05 14:19:33.554205"
    21 = "2021-10-05 14:19:33.554207"
    22 = "2021-10-05 14:19:33.554210"
    23 = "2021-10-05 14:19:33.554212"
    24 = "2021-10-05 14:19:33.554215"
    25 = "2021-10-05 14:19:33.554217"
    26 = "2021-10-05 14:19:33.554220"
    27 = "2021-10-05 14:19:33.554223"
    28 = "2021-10-05 14:19:33.554225"
    29 = "2021-10-05 14:19:33.554228"
    30 = "2021-10-05 14:19:33.554231"
    31 = "2021-10-05 14:19:33.554233"
    32 = "2021-10-05 14:19:33.554235"
    33 = "2021-10-05 14:19:33.554238"
    34 = "2021-10-05 14:19:33.554241"
    35 = "2021-10-05 14:19:33.554243"
    36 = "2021-10-05 14:19:33.554246"
    37 = "2021-10-05 14:19:33.554248"
    38 = "2021-10-05 14:19:33.554251"
    39 = "2021-10-05 14:19:33.554253"
    40 = "2021-10-05 14:19:33.554256"
    41 = "2021-10-05 14:19:33.554258"
    42 = "2021-10-05 14:19:33.554261"
    43 = "2021-10-05 14:19:33.554264"
    44 = "2021-10-05 14:19:33.554267"
    45 = "2021-10-05 14:19:33.554269"
    46 = "2021-10-05 14:19:33.554272"
    47 = "2021-10-05 14:19:33.554275"
    48 = "2021-10-05 14:19:33.554277"
    49 = "2021-10-05 14:19:33.554280"
    50 = "2021-10-05 14:19:33.554282"
    51 = "2021-10-05 14:19:33.554285"
    52 = "2021-10-05 14:19:33.554287"
    53 = "2021-10-05 14:19:33.554290"
    54 = "2021-10-05 14:19:33.554293"
    55 = "2021-10-05 14:19:33.554296"
    56 = "2021-10-05 14:19:33.554298"
    57 = "2021-10-05 14:19:33.554301"
    58 = "2021-10-05 14:19:33.554303"
    59 = "2021-10-05 14:19:33.554306"
    60 = "2021-10-05 14:19:33.554308"
    61 = "2021-10-05 14:19:33.554311"
    62 = "2021-10-05 14:19:33.554314"
    63 = "2021-10-05 14:19:33.554316"
    64 = "2021-10-05 14:19:33.554319"
    65 = "2021-10-05 14:19:33.554321"
    66 = "2021-10-05 14:19:33.554325"
    67 = "2021-10-05 14:19:33.554328"
    68 = "2021-10-05 14:19:33.554331"
    69 = "2021-10-05 14:19:33.554334"
    70 = "2021-10-05 14:19:33.554337"
    71 = "2021-10-05 14:19:33.554339"
    72 = "2021-10-05 14:19:33.554342"
    73 = "2021-10-05 14:19:33.554345"
    74 = "2021-10-05 14:19:33.554347"
    75 = "2021-10-05 14:19:33.554350"
    76 = "2021-10-05 14:19:33.554352"
    77 = "2021-10-05 14:19:33.554355"
    78 = "2021-10-05 14:19:33.554360"
    79 = "2021-10-05 14:19:33.554363"
    80 = "2021-10-05 14:19:33.554366"
    81 = "2021-10-05 14:19:33.554368"
    82 = "2021-10-05 14:19:33.554371"
    83 = "2021-10-05 14:19:33.554374"
    84 = "2021-10-05 14:19:33.554376"
    85 = "2021-10-05 14:19:33.554379"
    86 = "2021-10-05 14:19:33.554381"
}

foreach ($row in $timestamps.Keys) {
    $dataSheet.Cells.Item([int]$row, 6).Value = $timestamps[$row]
}

# --- Add the new "metadata" sheet after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

function Set-HeaderCell($sheet, $row, $col, $text) {
    $c = $sheet.Cells.Item($row, $col)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

Set-HeaderCell $metaSheet 1 2 "data_name"
Set-HeaderCell $metaSheet 1 3 "data_id"
Set-HeaderCell $metaSheet 1 4 "data_version"
Set-HeaderCell $metaSheet 1 5 "data_version_created"
Set-HeaderCell $metaSheet 1 6 "panel_query_time"
Set-HeaderCell $metaSheet 1 7 "panel_get_request"

# Row 2 - first cell (A2) uses same header-like style, rest plain
$a2 = $metaSheet.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$metaSheet.Cells.Item(2, 2).Value = "Childhood solid tumours cancer susceptibility"
$metaSheet.Cells.Item(2, 3).Value = 259

$d2 = $metaSheet.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1.15"
$d2.Style = "Normal"

$metaSheet.Cells.Item(2, 5).Value = "2021-07-15T09:49:01.721159Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:19:33.550853"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/259/?format=json"

Write-Output "edit complete"
